# correção nos dados e inicio da analise PNAD 2009
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 2 header labels: the "unnamed: 1_level_1" / "unnamed: 5_level_1"
#     placeholders are replaced by the real "total" label, and the
#     formerly-blank outer columns now also read "total".
$ws.Range("B2").Value = "total"
$ws.Range("C2").Value = "total"
$ws.Range("D2").Value = "condição de ocupação na semana de referência"
$ws.Range("E2").Value = "condição de ocupação na semana de referência.1"
$ws.Range("F2").Value = "total"

# --- Remove the now-redundant "situação do domicílio" sub-header row (5)
#     and the "grandes regiões" sub-header row (8, before the first delete
#     shifts rows up). Deleting the higher-numbered row first keeps the
#     original row numbering valid for the second delete.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()
